# Critical Reflection.docx edit:
# Append an extra sentence to the end of the "Difficulty assessment" paragraph,
# right after "...g each entry together with a space between entries."
$d = $word.ActiveDocument

$target = "g each entry together with a space between entries."
$addition = " Using a list instead of storing the entries as a string would have worked better for reversing such order."

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    if ($text -like "*$target*") {
        $r = $p.Range
        # Remember where the paragraph's text ends (Range.End is exclusive of
        # the paragraph mark, so inserting right before it appends a new run
        # after the existing "...entries." run without altering it).
        $insertStart = $r.End - 1
        $r.InsertAfter($addition)
        $insertEnd = $r.End - 1

        # Match the font formatting used throughout the paragraph so the new
        # run's rPr mirrors the existing "Franklin Gothic Book" runs.
        $newRange = $d.Range($insertStart, $insertEnd)
        $newRange.Font.Name = "Franklin Gothic Book"
        break
    }
}
